$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text updates (values that Excel will not mistake for numbers) ---
$ws.Range('D2').Value = '66.209.48'
$ws.Range('E2').Value = '  -0.12%  '
$ws.Range('D3').Value = '3.560.83'
$ws.Range('E3').Value = '  +1.42%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('E5').Value = '  -0.41%  '
$ws.Range('E6').Value = '  -1.12%  '
$ws.Range('D7').Value = '3.560.01'
$ws.Range('E7').Value = '  +1.40%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  +2.55%  '
$ws.Range('E10').Value = '  -0.40%  '
$ws.Range('E11').Value = '  -1.63%  '
$ws.Range('E12').Value = '  -0.33%  '
$ws.Range('D13').Value = '4.162.70'
$ws.Range('E13').Value = '  +1.41%  '
$ws.Range('E14').Value = '  -1.40%  '
$ws.Range('E15').Value = '  -0.93%  '
$ws.Range('D16').Value = '3.578.66'
$ws.Range('E16').Value = '  +1.99%  '
$ws.Range('D17').Value = '66.260.25'
$ws.Range('E17').Value = '  -0.07%  '
$ws.Range('E18').Value = '  -0.66%  '
$ws.Range('E19').Value = '  +8.56%  '
$ws.Range('E20').Value = '  -0.17%  '
$ws.Range('E21').Value = '  -1.02%  '
$ws.Range('E22').Value = '  +0.58%  '
$ws.Range('E23').Value = '  +1.63%  '
$ws.Range('E24').Value = '  +2.11%  '
$ws.Range('D25').Value = '3.702.17'
$ws.Range('E25').Value = '  +1.59%  '
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('E27').Value = '  -2.97%  '
$ws.Range('E28').Value = '  +0.70%  '
$ws.Range('E29').Value = '  -2.57%  '
$ws.Range('E30').Value = '  -1.80%  '
$ws.Range('E31').Value = '  +0.09%  '
$ws.Range('D32').Value = '3.555.99'
$ws.Range('E32').Value = '  +1.57%  '
$ws.Range('E33').Value = '  +0.58%  '
$ws.Range('E35').Value = '  -9.00%  '
$ws.Range('E36').Value = '  +0.05%  '
$ws.Range('E37').Value = '  -0.12%  '
$ws.Range('E38').Value = '  -1.45%  '
$ws.Range('E39').Value = '  -1.27%  '
$ws.Range('E40').Value = '  +1.91%  '
$ws.Range('E41').Value = '  -1.48%  '
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('E43').Value = '  +0.32%  '
$ws.Range('E44').Value = '  +0.83%  '
$ws.Range('E45').Value = '  +1.02%  '
$ws.Range('E46').Value = '  -0.03%  '
$ws.Range('E47').Value = '  -1.66%  '
$ws.Range('E49').Value = '  -1.81%  '
$ws.Range('E50').Value = '  -0.96%  '
$ws.Range('E51').Value = '  +1.74%  '

# --- Numeric-looking text updates (Price column values such as "0.490", "24.78") ---
# These must stay stored as text (matching the original inline-string cells),
# so force text format, assign, then restore the default "Normal" style so no
# visible formatting change is left behind.
$numericTextValues = @{
    'D5' = '604.93'
    'D6' = '143.51'
    'D9' = '0.490'
    'D10' = '0.136'
    'D14' = '0.0000206'
    'D15' = '30.07'
    'D19' = '11.35'
    'D20' = '6.17'
    'D21' = '14.71'
    'D22' = '428.90'
    'D23' = '0.608'
    'D24' = '79.67'
    'D29' = '9.04'
    'D30' = '7.85'
    'D31' = '0.999'
    'D33' = '25.38'
    'D37' = '7.78'
    'D39' = '5.54'
    'D40' = '173.79'
    'D41' = '0.0847'
    'D43' = '0.888'
    'D45' = '45.91'
    'D47' = '1.19'
    'D48' = '24.78'
    'D49' = '2.40'
    'D50' = '7.10'
    'D51' = '23.10'
}
foreach ($addr in $numericTextValues.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $numericTextValues[$addr]
    $cell.Style = "Normal"
}
